$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old rows 44-48 (generic Vol 1-5 entries) and rewrite rows 44-51 with detailed entries
$ws.Rows("44:48").ClearContents()

$ws.Cells.Item(44, 1).Value = 2003
$ws.Cells.Item(44, 2).Value = 'トーキョーN◎VA The Detonation スーパー・シナリオ・サポート Vol. 1 フロントライン'
$ws.Cells.Item(44, 3).Value = 'Tokyo Nova The Detonation Super Scenario Support Vol. 1: Froneline'
$ws.Cells.Item(44, 4).Value = 'Game Field'
$ws.Cells.Item(44, 5).Value = 'sssdet1a.jpg'
$ws.Cells.Item(44, 6).Value = 'scenario'

$ws.Cells.Item(45, 1).Value = 2004
$ws.Cells.Item(45, 2).Value = 'トーキョーN◎VA The Detonation スーパー・シナリオ・サポート Vol. 2 魂こがして'
$ws.Cells.Item(45, 3).Value = 'Tokyo Nova The Detonation Super Scenario Support Vol. 2: Flaming Fish'
$ws.Cells.Item(45, 4).Value = 'Game Field'
$ws.Cells.Item(45, 5).Value = 'sssdet2a.jpg'
$ws.Cells.Item(45, 6).Value = 'scenario'

$ws.Cells.Item(46, 1).Value = 2004
$ws.Cells.Item(46, 2).Value = 'トーキョーN◎VA The Detonation スーパー・シナリオ・サポート Vol. 3 フロントライン'
$ws.Cells.Item(46, 3).Value = 'Tokyo Nova The Detonation Super Scenario Support Vol. 3: Calling You'
$ws.Cells.Item(46, 4).Value = 'Game Field'
$ws.Cells.Item(46, 5).Value = 'sssdet3a.jpg'
$ws.Cells.Item(46, 6).Value = 'scenario'

$ws.Cells.Item(47, 1).Value = 2005
$ws.Cells.Item(47, 2).Value = 'トーキョーN◎VA The Detonation スーパー・シナリオ・サポート Vol. 1 この声がとどくまで'
$ws.Cells.Item(47, 3).Value = 'Tokyo Nova The Detonation Super Scenario Support Vol. 1: Crystal Wall'
$ws.Cells.Item(47, 4).Value = 'Game Field'
$ws.Cells.Item(47, 5).Value = 'sssdet1.jpg'
$ws.Cells.Item(47, 6).Value = 'scenario'

$ws.Cells.Item(48, 1).Value = 2005
$ws.Cells.Item(48, 2).Value = 'トーキョーN◎VA The Detonation スーパー・シナリオ・サポート Vol. 2 暗黒の世界'
$ws.Cells.Item(48, 3).Value = 'Tokyo Nova The Detonation Super Scenario Support Vol. 2: World of Darkness'
$ws.Cells.Item(48, 4).Value = 'Game Field'
$ws.Cells.Item(48, 5).Value = 'sssdet2.jpg'
$ws.Cells.Item(48, 6).Value = 'scenario'

$ws.Cells.Item(49, 1).Value = 2005
$ws.Cells.Item(49, 2).Value = 'トーキョーN◎VA The Detonation スーパー・シナリオ・サポート Vol. 3 仮面舞踏会'
$ws.Cells.Item(49, 3).Value = 'Tokyo Nova The Detonation Super Scenario Support Vol. 3: Behind the Mask'
$ws.Cells.Item(49, 4).Value = 'Game Field'
$ws.Cells.Item(49, 5).Value = 'sssdet3.jpg'
$ws.Cells.Item(49, 6).Value = 'scenario'

$ws.Cells.Item(50, 1).Value = 2005
$ws.Cells.Item(50, 2).Value = 'トーキョーN◎VA The Detonation スーパー・シナリオ・サポート Vol. 4 人間以上'
$ws.Cells.Item(50, 3).Value = 'Tokyo Nova The Detonation Super Scenario Support Vol. 4: Manplus'
$ws.Cells.Item(50, 4).Value = 'Game Field'
$ws.Cells.Item(50, 5).Value = 'sssdet4.jpg'
$ws.Cells.Item(50, 6).Value = 'scenario'

$ws.Cells.Item(51, 1).Value = 2004
$ws.Cells.Item(51, 2).Value = 'トーキョーN◎VA The Detonation スーパー・シナリオ・サポート Vol. 5 夕日の沈む朝'
$ws.Cells.Item(51, 3).Value = 'Tokyo Nova The Detonation Super Scenario Support Vol. 5: The Sun Against the Sun'
$ws.Cells.Item(51, 4).Value = 'Game Field'
$ws.Cells.Item(51, 5).Value = 'sssdet5.jpg'
$ws.Cells.Item(51, 6).Value = 'scenario'

$ws.Range("D44:D51").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1